$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (shifts all rows below it up by one)
$ws.Rows.Item(3).Delete()

# Update the selected cell / view to match the target state
$ws.Range("C20").Select()
